$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A27:AC27").ClearContents()
